{"js": "// Split the single run \"{m\" (and \"{m:\") into two separate runs \"{\" | \"m\"\n// (and \"{\" | \"m:\") without touching the text content, matching the\n// TokenIteratorFieldRewriterSplit parser change.\n//\n// Word (and this engine) re-tokenizes/merges adjacent same-format runs\n// whenever a paragraph's text is edited in place, so a naive delete+\n// re-insert collapses the two pieces back into one run. Inserting a\n// temporary bookmark at the desired split point first forces the run\n// boundary; deleting that bookmark afterwards (without touching the\n// surrounding text) leaves the run split in place instead of re-merging it.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nasync function splitRunAt(paragraph, searchText, bookmarkName) {\n  const matches = paragraph.search(searchText, { matchCase: true });\n  matches.load(\"items\");\n  await context.sync();\n  const match = matches.items[0];\n\n  // Collapsed range right after the first character of the match (i.e. the\n  // boundary between \"{\" and the rest, e.g. \"m\" or \"m:\").\n  const rest = match.search(searchText.substring(1), { matchCase: true });\n  rest.load(\"items\");\n  await context.sync();\n  const splitPoint = rest.items[0].getRange(\"Start\");\n\n  splitPoint.insertBookmark(bookmarkName);\n  await context.sync();\n\n  context.document.deleteBookmark(bookmarkName);\n  await context.sync();\n}\n\n// Paragraph \"name = {m:v.name},\" -> split \"{m\" into \"{\" + \"m\"\nawait splitRunAt(paragraphs.items[2], \"{m\", \"TempSplit1\");\n\n// Paragraph \"{m:endfor}\" -> split \"{m:\" into \"{\" + \"m:\"\nawait splitRunAt(paragraphs.items[4], \"{m:\", \"TempSplit2\");\n", "ps1": "# Split the single run \"{m\" (and \"{m:\") into two separate runs \"{\" | \"m\"\n# (and \"{\" | \"m:\") without touching the text content, matching the\n# TokenIteratorFieldRewriterSplit parser change.\n#\n# Trick: Word merges/\"re-tokenizes\" adjacent same-format runs whenever a\n# paragraph's text is edited in place. Inserting a temporary bookmark at the\n# desired split point first forces the run boundary; deleting that bookmark\n# afterwards (without touching the surrounding text) leaves the run split in\n# place instead of re-merging it.\n$d = $word.ActiveDocument\n\nfunction Split-RunAt($paragraphIndex, $searchText, $offsetIntoMatch, $bookmarkName) {\n    $pRange = $d.Paragraphs.Item($paragraphIndex).Range\n    $find = $pRange.Find\n    $find.ClearFormatting()\n    $find.Text = $searchText\n    $find.Execute() | Out-Null\n\n    $splitPos = $pRange.Start + $offsetIntoMatch\n    $splitRange = $d.Range($splitPos, $splitPos)\n    $d.Bookmarks.Add($bookmarkName, $splitRange) | Out-Null\n    $d.Bookmarks.Item($bookmarkName).Delete()\n}\n\n# Paragraph 3: \"name = {m:v.name},\" -> split \"{m\" into \"{\" + \"m\"\nSplit-RunAt 3 \"{m\" 1 \"TempSplit1\"\n\n# Paragraph 5: \"{m:endfor}\" -> split \"{m:\" into \"{\" + \"m:\"\nSplit-RunAt 5 \"{m:\" 1 \"TempSplit2\"\n"}
